$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3 V 0.3")

# Fix E5: convert from text to a real number
$ws.Cells.Item(5, 5).Value = 532900

# Add new row 6 with the breakout data
$ws.Cells.Item(6, 1).Value = "12/06/2024 09:45:47"
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = "PAISALO"
$ws.Cells.Item(6, 4).Value = "Paisalo Digital Ltd"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "532900"
$ws.Cells.Item(6, 6).Value = 8.56
$ws.Cells.Item(6, 7).Value = 69
$ws.Cells.Item(6, 8).Value = 4788944
